# AZ2522021|2:37PM Adding Parameters for E-1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimpleSearch")

# --- Fill in new parameter rows (A7:B12): URLs first (top to bottom), then
#     the matching labels back in (bottom to top) - the order the author typed them ---
$ws.Range("B7").Value  = "https://parts.z2data.com/RiskManager/Forecast?BomId=119090"
$ws.Range("B8").Value  = "https://parts.z2data.com/RiskManager/Compliance?BomId=119090"
$ws.Range("B9").Value  = "https://parts.z2data.com/RiskManager/Mitigation?BomId=119090"
$ws.Range("B10").Value = "https://parts.z2data.com/RiskManager/Report?BomId=119090"
$ws.Range("B11").Value = "https://parts.z2data.com/RiskManager/Scrub?BomId=119090"
$ws.Range("B12").Value = "https://parts.z2data.com/RiskManager?BomId=119090"

$ws.Range("A12").Value = "Pom_Dashboard_URL"
$ws.Range("A11").Value = "Scrub_Tab_URL"
$ws.Range("A10").Value = "Reports__Tab_URL"
$ws.Range("A9").Value  = "Mitigation_Tab_URL"
$ws.Range("A8").Value  = "Compliance_Tab_URL"
$ws.Range("A7").Value  = "ForeCast_Tab_URL"

# A few extra (still empty) formatted rows under the table, same as the author left them
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("B17").Value = ""

# --- Re-apply consistent direct formatting down the two columns ---
# Column A (labels) keeps the same look used for the existing rows
$ws.Range("A6").Copy()
$ws.Range("A7:A12").PasteSpecial(-4122)

# Column B (values) adopts the same look already used on B4 (border + vertical
# centered Calibri) instead of each row keeping its own one-off style
$ws.Range("B4").Copy()
$ws.Range("B5:B17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Resize the two columns now that they hold longer text ---
$ws.Columns("A:B").EntireColumn.AutoFit()

# --- Update selection to match where the author left off ---
$ws.Range("B3").Select()
